$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap date values between row 2 and row 4 (column D)
$ws.Range("D2").Value = 44971
$ws.Range("D4").Value = 44973

# Swap volume values between row 2 and row 4 (column M)
$ws.Range("M2").Value = 25
$ws.Range("M4").Value = 55
